{"js": "// The bug list's first five entries get replaced by two new \"flying fish\"\n// bug entries followed by a single blank paragraph (the two blank\n// paragraphs that already separated the list from \"Todo:\" are left alone).\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Paragraph 0: \"Fix the next turn button ...\" -> \"flying fish icons\"\nitems[0].insertText(\"flying fish icons\", Word.InsertLocation.replace);\n\n// Paragraph 1: \"Fix the \"Rewards\" text ...\" -> \"flying fish \u2013 with identifier bug?\"\nitems[1].insertText(\"flying fish \u2013 with identifier bug?\", Word.InsertLocation.replace);\n\n// Paragraph 2: \"Lone wolf doesn't lose modifier ???\" -> now a blank paragraph\nitems[2].insertText(\"\", Word.InsertLocation.replace);\n\n// Paragraphs 3 (\"Movement modifiers MoveSet at start of turn\") and 4\n// (\"Cauterize for modified max health lowers way more than It should\")\n// are removed entirely.\nitems[3].delete();\nitems[4].delete();\n\nawait context.sync();\n", "ps1": "# The bug list's first five entries get replaced by two new \"flying fish\"\n# bug entries followed by a single blank paragraph (the two blank\n# paragraphs that already separated the list from \"Todo:\" are left alone).\n$d = $word.ActiveDocument\n\nfunction Set-ParagraphText($para, [string]$text) {\n    # Paragraph.Range includes the trailing paragraph mark; trim it off so\n    # that assigning .Text rewrites every run's content in place without\n    # deleting the paragraph mark itself (works even when the paragraph has\n    # more than one run, e.g. \"Lone wolf doesn't lose modifier\" + \" ???\").\n    $r = $para.Range\n    $r.MoveEnd(1, -1) | Out-Null\n    $r.Text = $text\n}\n\n# Paragraph 1: \"Fix the next turn button ...\" -> \"flying fish icons\"\nSet-ParagraphText $d.Paragraphs.Item(1) \"flying fish icons\"\n\n# Paragraph 2: \"Fix the \"Rewards\" text ...\" -> \"flying fish \u2013 with identifier bug?\"\nSet-ParagraphText $d.Paragraphs.Item(2) \"flying fish \u2013 with identifier bug?\"\n\n# Paragraph 3: \"Lone wolf doesn't lose modifier ???\" -> now a blank paragraph\nSet-ParagraphText $d.Paragraphs.Item(3) \"\"\n\n# Paragraphs 5 (\"Cauterize for modified max health lowers way more than It\n# should\") and 4 (\"Movement modifiers MoveSet at start of turn\") are removed\n# entirely. Delete the higher index first so the lower index stays valid.\n$d.Paragraphs.Item(5).Range.Delete()\n$d.Paragraphs.Item(4).Range.Delete()\n"}
